# add data until March 8th
# Appends 5 new daily rows (2021-03-03 .. 2021-03-07) to the "Planilha1"
# sheet, continuing the existing daily COVID bulletin series (rows 2..244,
# dates 44015..44257) with rows 245..249 (dates 44258..44262).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$newRows = @(
    @(44258, 6063, 255, 1471, 7759, 1247, 191, 8, 183, 33),
    @(44259, 6075, 286, 1489, 7850, 1258, 198, 7, 191, 33),
    @(44260, 6132, 267, 1510, 7909, 1278, 197, 10, 187, 35),
    @(44261, 6174, 215, 1528, 7917, 1278, 214, 10, 204, 36),
    @(44262, 6223, 162, 1555, 7940, 1282, 237, 11, 226, 36)
)

$startRow = 245
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}

# Keep the sheet's view in sync with the data that was just appended,
# mirroring the selection move the author made after typing the new rows.
$null = $ws.Activate()
$lastRow = $startRow + $newRows.Count - 1
$null = $ws.Range("B" + ($lastRow + 1)).Select()
